$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.906.68"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.15%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.814.65"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.64%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.30"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4659"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3665"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07348"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8680"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.32"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.857.18"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.381"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07107"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.516"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.65"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.76%  "

$ws.Range("E17").Value = "  +0.24%  "

$ws.Range("E18").Value = "  +0.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.17%  "

$ws.Range("E20").Value = "  -0.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.940.07"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.300"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.64"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.082.76"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.893"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.01"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.37%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.30"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.139"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.260"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.28"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08898"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7552"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.156"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.42%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.488"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.913"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.82%  "

$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.085"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05275"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.51%  "

$ws.Range("E39").Value = "  -0.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.981"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.250"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5301"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.289"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.15%  "

$ws.Range("E44").Value = "  -0.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.427"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.95%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4871"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.86%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.39"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.18%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.26"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.659"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06291"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.09%  "
